$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for 2000, 2005, 2006, 2007, 2008, 2009 (old rows 2-7).
# This shifts the remaining rows (2010-2013, old rows 8-11) up to become
# rows 2-5.
$ws.Range("A2:G7").EntireRow.Delete()
